$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calibration")

# Solver-calibrated hazard-rate inputs (column A, rows 6-16). Setting these
# drives the downstream survival-probability / pricing formulas already
# present on the sheet to recalculate.
$ws.Range("A6").Value = 0.01764648405164415
$ws.Range("A7").Value = 0.02496482992018383
$ws.Range("A8").Value = 0.024991526205603725
$ws.Range("A9").Value = 0.027430204555949866
$ws.Range("A10").Value = 0.02783583605029908
$ws.Range("A11").Value = 0.03494769627555138
$ws.Range("A12").Value = 0.03494769627555141
$ws.Range("A13").Value = 0.03745852757770325
$ws.Range("A14").Value = 0.04061213575784079
$ws.Range("A15").Value = 0.0421056080273951
$ws.Range("A16").Value = 0.042115562055929645

# Update the sheet selection to match the saved view (A6 active cell,
# no longer scrolled to column V).
$ws.Activate()
$ws.Range("A6").Select()
